# Insert a new "MAE" column before the "Tipo" column, shifting Tipo to column E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing column D ("Tipo") to E
$ws.Range("D1").EntireColumn.Insert()

# New header and value for the MAE column
$ws.Range("D1").Value = "MAE"

$ws.Range("D2").Value = 0.1526989685211046
